$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells D2:E51 are stored as text (t="inlineStr") in the workbook, e.g. "322.03"
# and "-2.85%". A plain Range.Value assignment with a numeric-looking string
# makes Excel auto-convert the cell to a Number, which would change the cell's
# stored type. Force the target range to Text format first so the new values
# are written as literal text, then clear the temporary formatting so no
# extra/residual number-format styling is left behind on the cells.
$updateRange = $ws.Range("D2:E51")
$updateRange.NumberFormat = "@"

$ws.Range("D2").Value = "321.73"
$ws.Range("E2").Value = "-2.92%"
$ws.Range("D3").Value = "42.54"
$ws.Range("E3").Value = "-6.29%"
$ws.Range("D4").Value = "5.164"
$ws.Range("E4").Value = "-8.04%"
$ws.Range("D5").Value = "0.08191"
$ws.Range("E5").Value = "-1.91%"
$ws.Range("D6").Value = "4.289"
$ws.Range("D7").Value = "1.801"
$ws.Range("D8").Value = "0.9311"
$ws.Range("E8").Value = "-3.87%"
$ws.Range("D9").Value = "0.1109"
$ws.Range("E9").Value = "-5.90%"
$ws.Range("D10").Value = "0.1870"
$ws.Range("E10").Value = "-2.64%"
$ws.Range("D11").Value = "0.09499"
$ws.Range("E11").Value = "-3.69%"
$ws.Range("D12").Value = "0.04669"
$ws.Range("E12").Value = "1.11%"
$ws.Range("D13").Value = "7.427"
$ws.Range("E13").Value = "-27.98%"
$ws.Range("D14").Value = "0.1057"
$ws.Range("E14").Value = "-0.25%"
$ws.Range("D15").Value = "0.001288"
$ws.Range("E15").Value = "0.64%"
$ws.Range("D16").Value = "0.005727"
$ws.Range("E16").Value = "-5.27%"
$ws.Range("E17").Value = "-0.37%"
$ws.Range("D18").Value = "2.522"
$ws.Range("E18").Value = "-1.21%"
$ws.Range("D19").Value = "0.3374"
$ws.Range("E19").Value = "0.75%"
$ws.Range("D20").Value = "0.1388"
$ws.Range("E20").Value = "-0.37%"
$ws.Range("D21").Value = "0.2492"
$ws.Range("E21").Value = "-13.45%"
$ws.Range("D22").Value = "0.04158"
$ws.Range("E22").Value = "-0.76%"
$ws.Range("E23").Value = "-5.43%"
$ws.Range("D24").Value = "0.004427"
$ws.Range("E24").Value = "-2.98%"
$ws.Range("E25").Value = "-7.90%"
$ws.Range("D26").Value = "0.0002977"
$ws.Range("E26").Value = "-20.61%"
$ws.Range("D38").Value = "0.02774"
$ws.Range("D39").Value = "0.05597"
$ws.Range("E39").Value = "-2.96%"
$ws.Range("D40").Value = "0.008059"
$ws.Range("E40").Value = "2.29%"
$ws.Range("D41").Value = "0.1398"
$ws.Range("E41").Value = "-2.28%"
$ws.Range("D42").Value = "0.006542"
$ws.Range("E42").Value = "-9.99%"
$ws.Range("E43").Value = "2.93%"
$ws.Range("D44").Value = "0.008326"
$ws.Range("E44").Value = "-8.63%"
$ws.Range("D45").Value = "0.3505"
$ws.Range("E45").Value = "-1.18%"
$ws.Range("D46").Value = "0.00006951"
$ws.Range("E46").Value = "-2.68%"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").Value = "-0.33%"
$ws.Range("D48").Value = "0.003479"
$ws.Range("E48").Value = "-0.81%"
$ws.Range("D49").Value = "0.003527"
$ws.Range("E49").Value = "0.58%"
$ws.Range("D50").Value = "0.00002098"
$ws.Range("E50").Value = "-0.33%"
$ws.Range("D51").Value = "0.0001998"
$ws.Range("E51").Value = "-0.33%"

$updateRange.ClearFormats()
